$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.952.04"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "1.891.64"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.8304"
$ws.Range("E5").Value = "  +8.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.41"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3216"
$ws.Range("E8").Value = "  +5.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.59"
$ws.Range("E9").Value = "  +4.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07013"
$ws.Range("E10").Value = "  +2.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08035"
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7466"
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.902.82"
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.191"
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.25"
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").Value = "29.954.84"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.02"
$ws.Range("E17").Value = "  +1.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.914"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.16"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007745"
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "2.153.96"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.917"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1585"
$ws.Range("E25").Value = "  +23.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.72"
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.187"
$ws.Range("E27").Value = "  -0.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.82"
$ws.Range("E28").Value = "  +0.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.089"
$ws.Range("E29").Value = "  +2.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.374"
$ws.Range("E30").Value = "  -1.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.515"
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.250"
$ws.Range("E32").Value = "  -0.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05650"
$ws.Range("E33").Value = "  +7.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.068"
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.273"
$ws.Range("E35").Value = "  +1.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7305"
$ws.Range("E36").Value = "  +0.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.720"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01905"
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.771"
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4406"
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "71.82"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.943"
$ws.Range("E42").Value = "  -3.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8448"
$ws.Range("E43").Value = "  +1.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.887"
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.23"
$ws.Range("E46").Value = "  +1.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.592"
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.670"
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "988.38"
$ws.Range("E49").Value = "  +8.70%  "
$ws.Range("D50").Value = "2.050.66"
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.97"
$ws.Range("E51").Value = "  -0.36%  "
